# STM32: Make mask optional, reorder cmds, add GET_DESTINATION_BUFFER cmd
#
# UART Commands sheet, rows 28/29 currently hold (in this order):
#   28: 0x62 SET_DESTINATION_BUFFER
#   29: 0x63 UPDATE_SCROLL_BUFFER
# They get reordered so UPDATE_SCROLL_BUFFER comes first, then
# SET_DESTINATION_BUFFER - the 0x62/0x63 codes in column C stay where they
# are, only the Name/Payload-Length/Description columns move.
# Three brand new commands are appended afterwards:
#   30: 0x64 GET_DESTINATION_BUFFER
#   31: 0x65 SET_MASK_ENABLED
#   32: 0x66 GET_MASK_ENABLED

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UART Commands")

# --- 1) make room for three new rows right after row 29 -------------------
$ws.Range("A30:A32").EntireRow.Insert()

# --- 2) swap the Name / Payload Length / Description of rows 28 and 29 ----
$d28 = $ws.Range("D28").Value2
$e28 = $ws.Range("E28").Value2
$g28 = $ws.Range("G28").Value2

$ws.Range("D28").Value = $ws.Range("D29").Value2
$ws.Range("E28").Value = $ws.Range("E29").Value2
$ws.Range("G28").Value = $ws.Range("G29").Value2

$ws.Range("D29").Value = $d28
$ws.Range("E29").Value = $e28
$ws.Range("G29").Value = $g28

# row heights for the swapped rows (UPDATE_SCROLL_BUFFER's description is
# long, SET_DESTINATION_BUFFER's is short)
$ws.Range("A28").RowHeight = 405
$ws.Range("A29").RowHeight = 75

# --- 3) fill in the three new commands -------------------------------------
# (text cells are written in the same order the original author entered
# them in, so new shared-string entries line up the same way)

# 0x64 GET_DESTINATION_BUFFER
$ws.Range("C30").Value = "0x64"
$ws.Range("D30").Value = "GET_DESTINATION_BUFFER"
$ws.Range("G30").Value = "Get the SPI destination buffer.`nResponse:`nByte 0 - The buffer ID as expected by SET_DESTINATION_BUFFER"
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 1

# Names for the two mask rows first ...
$ws.Range("D31").Value = "SET_MASK_ENABLED"
$ws.Range("D32").Value = "GET_MASK_ENABLED"

# ... then their command codes ...
$ws.Range("C31").Value = "0x65"
$ws.Range("C32").Value = "0x66"

# ... then their descriptions.
$ws.Range("G31").Value = "Set mask compositing on or off.`nParameters:`nByte 0 - 1 for on, 0 for off"
$ws.Range("G32").Value = "Return mask compositing state.`nResponse:`nByte 0 - 1 for on, 0 for off"

# 0x65 SET_MASK_ENABLED
$ws.Range("E31").Value = 1
$ws.Range("F31").Value = 0

# 0x66 GET_MASK_ENABLED
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 1

$ws.Range("A30").RowHeight = 45
$ws.Range("A31").RowHeight = 45
$ws.Range("A32").RowHeight = 45

# --- 4) leave the selection where the author ended up ----------------------
$ws.Range("C32").Select()
